$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metricas")

# Fill in row 9 ("Implementar PilaCL") with its metrics, mirroring rows 6-8.
$ws.Range("A9").Value = "Implementar PilaCL"
$ws.Range("B9").Value = 25
$ws.Range("C9").Value = 22
$ws.Range("D9").Value = 0.00694444444444444406
$ws.Range("E9").Value = 0.83333333333333337034
$ws.Range("F9").Value = 0.84097222222222223209
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0

# G9 and J9 already hold formulas (=F9-E9 and =G9+I9 respectively); they
# will recalculate automatically, but re-assert them to be safe.
$ws.Range("G9").Formula = "=F9-E9"
$ws.Range("J9").Formula = "=G9+I9"

# Move the active selection, as recorded in the saved view state.
$ws.Range("A9").Select()

$wb.Save()
